$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.981.01"
$ws.Range("E2").Value = "  +9.80%  "

$ws.Range("D3").Value = "3.464.15"
$ws.Range("E3").Value = "  +6.04%  "

$ws.Range("E4").Value = "  +0.27%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "414.28"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.92%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "122.50"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +12.52%  "

$ws.Range("D7").Value = "3.455.93"
$ws.Range("E7").Value = "  +5.96%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.594"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.71%  "

$ws.Range("E9").Value = "  +0.13%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.684"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +10.53%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.130"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +35.61%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "41.20"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +5.03%  "

$ws.Range("E13").Value = "  +0.41%  "

$ws.Range("D14").Value = "4.023.25"
$ws.Range("E14").Value = "  +6.42%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "8.59"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +4.33%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "19.98"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +5.42%  "

$ws.Range("D17").Value = "3.465.62"
$ws.Range("E17").Value = "  +6.20%  "

$ws.Range("D18").Value = "62.932.79"
$ws.Range("E18").Value = "  +10.08%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.03"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "10.88"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.22%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.0000137"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +27.83%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "3.30"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "316.08"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +6.48%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "81.61"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +10.07%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "12.84"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.79%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.17"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "30.83"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +9.57%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.75"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +4.23%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.81"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.61%  "

$ws.Range("B30").Value = "LEO"
$ws.Range("C30").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.30"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.88%  "

$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.174"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.90%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.116"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +3.80%  "

$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.59"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +20.90%  "

$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "11.59"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.52%  "

$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.58%  "

$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "41.99"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +4.57%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0490"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.49%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "52.24"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.63%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.49"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.60%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.03"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.92%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.99"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +6.10%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.125"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +3.40%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "135.76"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.78%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.281"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.57%  "

$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "16.79"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.88"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.02%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.26"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.20%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "21.95"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.95%  "

$ws.Range("D50").Value = "2.180.61"
$ws.Range("E50").Value = "  +1.19%  "

$ws.Range("E51").Value = "  -0.01%  "
